# mode/LoneWolf/NewEngineFormulas.xlsx — "new engines and changes to economy"
#
# The only real input change is D2 (the "Last Level" HP target) going from
# 50 to 100. Every H20:H168 cell is a formula
#   ROUND(((F-$C$4)/($D$4-$C$4))*($D$2-$C$2),0)+$C$2
# so changing D2 ripples through automatically on recalculation — no need
# to touch those cells by hand.
#
# The sheet view also scrolled down (topLeftCell A133) and the selection
# moved to H9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Make sure we're working on/viewing the right sheet before touching the
# selection / scroll position.
$ws.Activate()

# Core data edit: bump the "Last Level" HP value used by every H-column
# formula. This single write causes the whole H20:H168 recalculation seen
# in the diff.
$ws.Range("D2").Value = 100

# View state: scroll so row 133 is the first visible row, then move the
# active selection to H9 (matches the sheetView/selection in the diff).
$excel.ActiveWindow.ScrollRow = 133
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H9").Select()
